# Refresh the crypto price table with the latest scrape (GitHub Actions run).
# Columns: A=rank(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of "Price" values are plain digits (e.g. "1.00", "6.41"); the sheet
# always stores this column as text (prices use "." as a thousands separator too,
# e.g. "57.974.83"), so a leading apostrophe is used to stop Excel from silently
# re-typing those cells as numbers.

$ws.Range("D2").Value = '57.974.83'
$ws.Range("E2").Value = '  -2.27%  '

$ws.Range("D3").Value = '2.570.42'
$ws.Range("E3").Value = '  -2.68%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").Value = '''517.31'
$ws.Range("E5").Value = '  -2.16%  '

$ws.Range("D6").Value = '''139.16'
$ws.Range("E6").Value = '  -4.10%  '

$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '  -0.05%  '

$ws.Range("E8").Value = '  -1.68%  '

$ws.Range("D9").Value = '2.582.31'
$ws.Range("E9").Value = '  -2.71%  '

$ws.Range("D10").Value = '''6.41'
$ws.Range("E10").Value = '  -3.90%  '

$ws.Range("D11").Value = '''0.0992'
$ws.Range("E11").Value = '  -5.11%  '

$ws.Range("D12").Value = '''0.326'
$ws.Range("E12").Value = '  -3.76%  '

$ws.Range("E13").Value = '  +0.44%  '

$ws.Range("D14").Value = '3.019.30'
$ws.Range("E14").Value = '  -2.80%  '

$ws.Range("D15").Value = '57.926.23'
$ws.Range("E15").Value = '  -2.26%  '

$ws.Range("D16").Value = '''20.06'
$ws.Range("E16").Value = '  -4.69%  '

$ws.Range("D17").Value = '2.572.37'
$ws.Range("E17").Value = '  -3.99%  '

$ws.Range("E18").Value = '  -4.37%  '

$ws.Range("D19").Value = '''332.55'
$ws.Range("E19").Value = '  -2.78%  '

$ws.Range("D20").Value = '''4.28'
$ws.Range("E20").Value = '  -4.23%  '

$ws.Range("D21").Value = '''10.06'
$ws.Range("E21").Value = '  -5.35%  '

$ws.Range("D22").Value = '''6.32'
$ws.Range("E22").Value = '  -0.46%  '

$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '  +0.02%  '

$ws.Range("D24").Value = '''65.77'
$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("E25").Value = '  -1.50%  '

$ws.Range("E26").Value = '  +0.10%  '

$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").Value = '''0.398'
$ws.Range("E27").Value = '  -5.25%  '

$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '2.690.89'
$ws.Range("E28").Value = '  -2.45%  '

$ws.Range("D29").Value = '''6.91'
$ws.Range("E29").Value = '  -4.85%  '

$ws.Range("D30").Value = '''0.998'
$ws.Range("E30").Value = '  -0.04%  '

$ws.Range("E31").Value = '  -10.50%  '

$ws.Range("E32").Value = '  -7.94%  '

$ws.Range("E33").Value = '  -3.83%  '

$ws.Range("D34").Value = '''18.61'
$ws.Range("E34").Value = '  -1.95%  '

$ws.Range("D35").Value = '''149.01'
$ws.Range("E35").Value = '  -0.82%  '

$ws.Range("E36").Value = '  -7.41%  '

$ws.Range("E37").Value = '  -7.92%  '

$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").Value = '''36.15'
$ws.Range("E38").Value = '  -1.25%  '

$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").Value = '''0.840'
$ws.Range("E39").Value = '  -3.32%  '

$ws.Range("D40").Value = '''0.823'
$ws.Range("E40").Value = '  -8.90%  '

$ws.Range("D41").Value = '''1.43'
$ws.Range("E41").Value = '  -5.04%  '

$ws.Range("D42").Value = '''3.47'
$ws.Range("E42").Value = '  -5.39%  '

$ws.Range("E43").Value = '  +0.02%  '

$ws.Range("D44").Value = '''273.52'
$ws.Range("E44").Value = '  +0.13%  '

$ws.Range("D45").Value = '''10.69'
$ws.Range("E45").Value = '  +0.29%  '

$ws.Range("D46").Value = '''0.588'
$ws.Range("E46").Value = '  -2.52%  '

$ws.Range("D47").Value = '''0.0939'
$ws.Range("E47").Value = '  -3.76%  '

$ws.Range("D48").Value = '''0.0514'
$ws.Range("E48").Value = '  -4.74%  '

$ws.Range("D49").Value = '''18.34'
$ws.Range("E49").Value = '  -5.74%  '

$ws.Range("D50").Value = '1.959.80'
$ws.Range("E50").Value = '  -4.02%  '

$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = '''4.49'
$ws.Range("E51").Value = '  -6.12%  '

# Entering those values with a leading apostrophe leaves a "quote prefix" style on
# the cell (Excel's normal way of remembering "this number-looking text was typed
# as text"). Re-apply the plain, unstyled look the rest of the table uses so those
# cells stay visually identical to their neighbours.
$plainStyle = $ws.Range("B2").Style
$numericLikeCells = "D5","D6","D7","D10","D11","D12","D16","D19","D20","D21","D22","D23","D24","D27","D29","D30","D34","D35","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49","D51"
foreach ($cellRef in $numericLikeCells) {
    $ws.Range($cellRef).Style = $plainStyle
}
